$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-04 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-05 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("34÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷6=", 2) | Out-Null
$d.Content.Find.Execute("86÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷8=", 2) | Out-Null
$d.Content.Find.Execute("50÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷8=", 2) | Out-Null
$d.Content.Find.Execute("10÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷3=", 2) | Out-Null
$d.Content.Find.Execute("85÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷4=", 2) | Out-Null
$d.Content.Find.Execute("45÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷8=", 2) | Out-Null
$d.Content.Find.Execute("21÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷6=", 2) | Out-Null
$d.Content.Find.Execute("27÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷5=", 2) | Out-Null
$d.Content.Find.Execute("75÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷2=", 2) | Out-Null
$d.Content.Find.Execute("53÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷8=", 2) | Out-Null
$d.Content.Find.Execute("49÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷4=", 2) | Out-Null
$d.Content.Find.Execute("39÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "96÷9=", 2) | Out-Null
$d.Content.Find.Execute("93÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷4=", 2) | Out-Null
$d.Content.Find.Execute("40÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷9=", 2) | Out-Null
$d.Content.Find.Execute("49÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷5=", 2) | Out-Null
$d.Content.Find.Execute("89÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷6=", 2) | Out-Null
$d.Content.Find.Execute("98÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷8=", 2) | Out-Null
$d.Content.Find.Execute("85÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "45÷4=", 2) | Out-Null
$d.Content.Find.Execute("36÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷5=", 2) | Out-Null
$d.Content.Find.Execute("61÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "73÷2=", 2) | Out-Null
$d.Content.Find.Execute("28÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷4=", 2) | Out-Null
$d.Content.Find.Execute("54÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷5=", 2) | Out-Null
$d.Content.Find.Execute("45÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷2=", 2) | Out-Null
$d.Content.Find.Execute("56÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷3=", 2) | Out-Null
$d.Content.Find.Execute("60÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷7=", 2) | Out-Null
